$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2 (Target cluster = ECs) ---
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.3987123333333333
$ws.Range("H2").Value = 1.196137
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 43.28121633333333
$ws.Range("N2").Value = 129.843649
$ws.Range("O2").Value = 0.533749049291363
$ws.Range("P2").Value = 0.533749049291363
$ws.Range("Q2").Value = 17.25675475376811
$ws.Range("R2").Value = 155.310792783913
$ws.Range("S2").Value = 0.533749049291363
$ws.Range("T2").Value = 0.533749049291363

# --- Update row 3 (Target cluster = FAPs) ---
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.3987123333333333
$ws.Range("H3").Value = 1.196137
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 9.865038
$ws.Range("N3").Value = 29.595114
$ws.Range("O3").Value = 0.1216568086527629
$ws.Range("P3").Value = 0.1216568086527629
$ws.Range("Q3").Value = 3.933312319402
$ws.Range("R3").Value = 35.399810874618
$ws.Range("S3").Value = 0.1216568086527629
$ws.Range("T3").Value = 0.1216568086527629

# --- Row 4 becomes a new record (Target cluster = M2), replacing old row4 data ---
$ws.Range("A4").Value = "M2"
$ws.Range("B4").Value = "Matn1"
$ws.Range("C4").Value = "Itga1"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.3987123333333333
$ws.Range("H4").Value = 1.196137
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 1.460025
$ws.Range("N4").Value = 4.380075
$ws.Range("O4").Value = 0.01800519998536753
$ws.Range("P4").Value = 0.01800519998536753
$ws.Range("Q4").Value = 0.5821299744749999
$ws.Range("R4").Value = 5.239169770275
$ws.Range("S4").Value = 0.01800519998536753
$ws.Range("T4").Value = 0.01800519998536753

# --- New row 5 (Target cluster = sCs), with updated values ---
$ws.Range("A5").Value = "M2"
$ws.Range("B5").Value = "Matn1"
$ws.Range("C5").Value = "Itga1"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.3987123333333333
$ws.Range("H5").Value = 1.196137
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 26.482795
$ws.Range("N5").Value = 79.448385
$ws.Range("O5").Value = 0.3265889420705065
$ws.Range("P5").Value = 0.3265889420705065
$ws.Range("Q5").Value = 10.55901698763833
$ws.Range("R5").Value = 95.031152888745
$ws.Range("S5").Value = 0.3265889420705065
$ws.Range("T5").Value = 0.3265889420705065
